$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Recorded-By column: reorder "dnasr281@gmail.com, System" so System comes
#    first -> "System, dnasr281@gmail.com" across every session row that had
#    both recorders.
# ---------------------------------------------------------------------------
$recordedByCells = @(
    "G8","G9","G10",
    "G34","G35","G36",
    "G60","G61","G62",
    "G86","G87","G88",
    "G112","G113","G114",
    "G138","G139","G140",
    "G164","G167",
    "G191","G194",
    "G218","G221",
    "G245","G248",
    "G272","G275",
    "G299","G302"
)

foreach ($addr in $recordedByCells) {
    $cell = $ws.Range($addr)
    $val = $cell.Value()
    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}

# ---------------------------------------------------------------------------
# 2) Daily attendance processing: session 14 (21/12/2025) got recorded, so it
#    swaps places (attendance values + Recorded/Not-Recorded row formatting)
#    with the still-pending session 13 (20/12/2025) row directly above it,
#    for every group block. Session number/date columns (A-F) stay put; only
#    the Recorded-By / Students / Status cells (G:I) and the row's
#    conditional-format colour move between the two rows.
# ---------------------------------------------------------------------------
$swapRowPairs = @(
    @{ Pending = 13;  Recorded = 14  },
    @{ Pending = 39;  Recorded = 40  },
    @{ Pending = 65;  Recorded = 66  },
    @{ Pending = 91;  Recorded = 92  },
    @{ Pending = 117; Recorded = 118 },
    @{ Pending = 143; Recorded = 144 }
)

$scratchRow = 5000

foreach ($pair in $swapRowPairs) {
    $r1 = $pair.Pending
    $r2 = $pair.Recorded

    $range1 = "A" + $r1 + ":I" + $r1
    $range2 = "A" + $r2 + ":I" + $r2
    $scratchRange = "A" + $scratchRow + ":I" + $scratchRow

    # Remember the G/H/I (Recorded By / Students / Status) values of both rows.
    $g1 = $ws.Range("G" + $r1).Value()
    $h1 = $ws.Range("H" + $r1).Value()
    $i1 = $ws.Range("I" + $r1).Value()
    $g2 = $ws.Range("G" + $r2).Value()
    $h2 = $ws.Range("H" + $r2).Value()
    $i2 = $ws.Range("I" + $r2).Value()

    # Swap the row-level formatting (fill colour etc.) between the two rows
    # using a scratch row as a temporary holding area.
    $ws.Range($range2).Copy()
    $ws.Range($scratchRange).PasteSpecial(-4122)
    $ws.Range($range1).Copy()
    $ws.Range($range2).PasteSpecial(-4122)
    $ws.Range($scratchRange).Copy()
    $ws.Range($range1).PasteSpecial(-4122)
    $ws.Range($scratchRange).Clear()

    # Restore the Recorded-By / Students / Status values, swapped.
    $ws.Range("G" + $r1).Value = $g2
    $ws.Range("H" + $r1).Value = $h2
    $ws.Range("I" + $r1).Value = $i2
    $ws.Range("G" + $r2).Value = $g1
    $ws.Range("H" + $r2).Value = $h1
    $ws.Range("I" + $r2).Value = $i1
}
